# Version 3.3 Incluido -R
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add three new TEST rows (rows 89-91) with TEMA number, PREGUNTAS topic and TOTAL count
$ws.Range("A89").Value = 16
$ws.Range("B89").Value = "LOPJ"
$ws.Range("C89").Value = 80

$ws.Range("A90").Value = 19
$ws.Range("B90").Value = "LOPJ"
$ws.Range("C90").Value = 80

$ws.Range("A91").Value = 24
$ws.Range("B91").Value = "GOBIERNO"
$ws.Range("C91").Value = 75

$excel.CalculateFull()
